# Wario - Level 6 begin
# Insert a new row for "Leave Level 5" right before the "Level 6" section header,
# fill in the new timing data, and extend the existing timing figures for the
# rows leading up to it (rows 45-49 gain a "Level entry" timestamp in column B,
# which in turn makes their column-D "Diff" formula produce a real value).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert the new row 50 ("Leave Level 5") -------------------------------
# This shifts every row from the old row 50 ("Level 6" header) onward down by
# one, which is exactly the shape of the target diff (old r -> new r+1).
$ws.Rows(50).Insert()

# Copy formatting (style s="4") from the row directly above so the new row
# matches the look of the other data rows, then overwrite with real data.
$ws.Range("A49:D49").Copy()
$ws.Range("A50:D50").PasteSpecial(-4122)

$ws.Range("A50").Value = "Leave Level 5"
$ws.Range("B50").Value = 17309
$ws.Range("C50").Value = 18801
$ws.Range("D50").Formula = "=IF(B50 >  0,C50-B50, 0)"

# --- Fill in the new "entry" timestamps for rows 45-49 ---------------------
$ws.Range("B45").Value = 15841
$ws.Range("B46").Value = 16245
$ws.Range("B47").Value = 16760
$ws.Range("B48").Value = 16939
$ws.Range("B49").Value = 17105

# --- The row right after the inserted row ("Enter Level 6") also gains a
# start timestamp (old row 51 -> new row 52) ---------------------------------
$ws.Range("B52").Value = 17558

# --- Restore the selection/active cell (view moved because of the new row) -
$ws.Range("B53").Select()
